# Auto-generated edit script applying numeric corrections to the
# Sheets workbook (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR), per the scheduled
# market-data refresh. Updates currentAveragePrice* / Leve* columns
# (H:N) for the affected leve rows on each sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(100, 8).Value = 1864.619  # H100
$ws.Cells.Item(100, 9).Value = 1727.2858  # I100
$ws.Cells.Item(100, 10).Value = 1933.2858  # J100
$ws.Cells.Item(100, 11).Value = 1727.2858  # K100
$ws.Cells.Item(100, 12).Value = 1933.2858  # L100
$ws.Cells.Item(100, 13).Value = -1186.2858  # M100
$ws.Cells.Item(100, 14).Value = -3015.2858  # N100

$ws.Cells.Item(137, 8).Value = 3174.1428  # H137
$ws.Cells.Item(137, 9).Value = 2774.8  # I137
$ws.Cells.Item(137, 10).Value = 4172.5  # J137
$ws.Cells.Item(137, 11).Value = 8324.400000000001  # K137
$ws.Cells.Item(137, 12).Value = 12517.5  # L137
$ws.Cells.Item(137, 13).Value = -5774.400000000001  # M137
$ws.Cells.Item(137, 14).Value = -17617.5  # N137

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1429723.5  # H2
$ws.Cells.Item(2, 9).Value = 2042153.4  # I2
$ws.Cells.Item(2, 10).Value = 720.6667  # J2
$ws.Cells.Item(2, 11).Value = 2042153.4  # K2
$ws.Cells.Item(2, 12).Value = 720.6667  # L2
$ws.Cells.Item(2, 13).Value = -2042040.4  # M2
$ws.Cells.Item(2, 14).Value = -946.6667  # N2

$ws.Cells.Item(34, 8).Value = 0  # H34
$ws.Cells.Item(34, 9).Value = 0  # I34
$ws.Cells.Item(34, 11).Value = 0  # K34
$ws.Cells.Item(34, 13).ClearContents()  # M34

$ws.Cells.Item(37, 8).Value = 7016.3335  # H37
$ws.Cells.Item(37, 9).Value = 2619.6  # I37
$ws.Cells.Item(37, 10).Value = 29000  # J37
$ws.Cells.Item(37, 11).Value = 2619.6  # K37
$ws.Cells.Item(37, 12).Value = 29000  # L37
$ws.Cells.Item(37, 13).Value = -2346.6  # M37
$ws.Cells.Item(37, 14).Value = -29546  # N37

$ws.Cells.Item(45, 8).Value = 9627.454  # H45
$ws.Cells.Item(45, 9).Value = 10374.235  # I45
$ws.Cells.Item(45, 11).Value = 10374.235  # K45
$ws.Cells.Item(45, 13).Value = -9997.235000000001  # M45

$ws.Cells.Item(55, 8).Value = 24799.2  # H55
$ws.Cells.Item(55, 10).Value = 24799.2  # J55
$ws.Cells.Item(55, 12).Value = 24799.2  # L55
$ws.Cells.Item(55, 14).Value = -25429.2  # N55

$ws.Cells.Item(102, 8).Value = 4618.3184  # H102
$ws.Cells.Item(102, 9).Value = 3357.3125  # I102
$ws.Cells.Item(102, 11).Value = 3357.3125  # K102
$ws.Cells.Item(102, 13).Value = -1735.3125  # M102

$ws.Cells.Item(110, 8).Value = 2960.4443  # H110
$ws.Cells.Item(110, 9).Value = 2684.2273  # I110
$ws.Cells.Item(110, 11).Value = 2684.2273  # K110
$ws.Cells.Item(110, 13).Value = -639.2273  # M110

$ws.Cells.Item(116, 8).Value = 1429723.5  # H116
$ws.Cells.Item(116, 9).Value = 2042153.4  # I116
$ws.Cells.Item(116, 10).Value = 720.6667  # J116
$ws.Cells.Item(116, 11).Value = 2042153.4  # K116
$ws.Cells.Item(116, 12).Value = 720.6667  # L116
$ws.Cells.Item(116, 13).Value = -2039859.4  # M116
$ws.Cells.Item(116, 14).Value = -5308.6667  # N116

$ws.Cells.Item(122, 8).Value = 2316.5557  # H122
$ws.Cells.Item(122, 9).Value = 1864.5883  # I122
$ws.Cells.Item(122, 10).Value = 10000  # J122
$ws.Cells.Item(122, 11).Value = 5593.7649  # K122
$ws.Cells.Item(122, 12).Value = 30000  # L122
$ws.Cells.Item(122, 13).Value = -3143.7649  # M122
$ws.Cells.Item(122, 14).Value = -34900  # N122

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1429723.5  # H3
$ws.Cells.Item(3, 9).Value = 2042153.4  # I3
$ws.Cells.Item(3, 10).Value = 720.6667  # J3
$ws.Cells.Item(3, 11).Value = 2042153.4  # K3
$ws.Cells.Item(3, 12).Value = 720.6667  # L3
$ws.Cells.Item(3, 13).Value = -2042039.4  # M3
$ws.Cells.Item(3, 14).Value = -948.6667  # N3

$ws.Cells.Item(105, 8).Value = 3850.7778  # H105
$ws.Cells.Item(105, 9).Value = 2300.7  # I105
$ws.Cells.Item(105, 10).Value = 8279.571  # J105
$ws.Cells.Item(105, 11).Value = 2300.7  # K105
$ws.Cells.Item(105, 12).Value = 8279.571  # L105
$ws.Cells.Item(105, 13).Value = -553.6999999999998  # M105
$ws.Cells.Item(105, 14).Value = -11773.571  # N105

$ws.Cells.Item(107, 8).Value = 5833.3335  # H107
$ws.Cells.Item(107, 9).Value = 4327.4546  # I107
$ws.Cells.Item(107, 11).Value = 4327.4546  # K107
$ws.Cells.Item(107, 13).Value = -2407.4546  # M107

$ws.Cells.Item(130, 8).Value = 84999.5  # H130
$ws.Cells.Item(130, 10).Value = 84999.5  # J130
$ws.Cells.Item(130, 12).Value = 84999.5  # L130
$ws.Cells.Item(130, 14).Value = -95039.5  # N130

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 3288.5454  # H16
$ws.Cells.Item(16, 9).Value = 2187.375  # I16
$ws.Cells.Item(16, 10).Value = 6225  # J16
$ws.Cells.Item(16, 11).Value = 2187.375  # K16
$ws.Cells.Item(16, 12).Value = 6225  # L16
$ws.Cells.Item(16, 13).Value = -1900.375  # M16
$ws.Cells.Item(16, 14).Value = -6799  # N16

$ws.Cells.Item(22, 8).Value = 1042.6786  # H22
$ws.Cells.Item(22, 9).Value = 874.0714  # I22
$ws.Cells.Item(22, 10).Value = 1211.2858  # J22
$ws.Cells.Item(22, 11).Value = 874.0714  # K22
$ws.Cells.Item(22, 12).Value = 1211.2858  # L22
$ws.Cells.Item(22, 13).Value = -524.0714  # M22
$ws.Cells.Item(22, 14).Value = -1911.2858  # N22

$ws.Cells.Item(31, 8).Value = 2562.5334  # H31
$ws.Cells.Item(31, 9).Value = 2562.5334  # I31
$ws.Cells.Item(31, 10).Value = 0  # J31
$ws.Cells.Item(31, 11).Value = 2562.5334  # K31
$ws.Cells.Item(31, 12).Value = 0  # L31
$ws.Cells.Item(31, 13).Value = -2267.5334  # M31
$ws.Cells.Item(31, 14).ClearContents()  # N31

$ws.Cells.Item(34, 8).Value = 2562.5334  # H34
$ws.Cells.Item(34, 9).Value = 2562.5334  # I34
$ws.Cells.Item(34, 10).Value = 0  # J34
$ws.Cells.Item(34, 11).Value = 2562.5334  # K34
$ws.Cells.Item(34, 12).Value = 0  # L34
$ws.Cells.Item(34, 13).Value = -2360.5334  # M34
$ws.Cells.Item(34, 14).ClearContents()  # N34

$ws.Cells.Item(58, 8).Value = 6913.5806  # H58
$ws.Cells.Item(58, 9).Value = 7511.1333  # I58
$ws.Cells.Item(58, 10).Value = 6353.375  # J58
$ws.Cells.Item(58, 11).Value = 7511.1333  # K58
$ws.Cells.Item(58, 12).Value = 6353.375  # L58
$ws.Cells.Item(58, 13).Value = -7308.1333  # M58
$ws.Cells.Item(58, 14).Value = -6759.375  # N58

$ws.Cells.Item(99, 8).Value = 12828.8  # H99
$ws.Cells.Item(99, 9).Value = 10483.417  # I99
$ws.Cells.Item(99, 10).Value = 14993.77  # J99
$ws.Cells.Item(99, 11).Value = 10483.417  # K99
$ws.Cells.Item(99, 12).Value = 14993.77  # L99
$ws.Cells.Item(99, 13).Value = -8985.416999999999  # M99
$ws.Cells.Item(99, 14).Value = -17989.77  # N99

$ws.Cells.Item(107, 8).Value = 830.5625  # H107
$ws.Cells.Item(107, 9).Value = 224  # I107
$ws.Cells.Item(107, 10).Value = 2165  # J107
$ws.Cells.Item(107, 11).Value = 224  # K107
$ws.Cells.Item(107, 12).Value = 2165  # L107
$ws.Cells.Item(107, 13).Value = 1696  # M107
$ws.Cells.Item(107, 14).Value = -6005  # N107

$ws.Cells.Item(113, 8).Value = 3288.5454  # H113
$ws.Cells.Item(113, 9).Value = 2187.375  # I113
$ws.Cells.Item(113, 10).Value = 6225  # J113
$ws.Cells.Item(113, 11).Value = 2187.375  # K113
$ws.Cells.Item(113, 12).Value = 6225  # L113
$ws.Cells.Item(113, 13).Value = -17.375  # M113
$ws.Cells.Item(113, 14).Value = -10565  # N113

$ws.Cells.Item(126, 8).Value = 12828.8  # H126
$ws.Cells.Item(126, 9).Value = 10483.417  # I126
$ws.Cells.Item(126, 10).Value = 14993.77  # J126
$ws.Cells.Item(126, 11).Value = 31450.251  # K126
$ws.Cells.Item(126, 12).Value = 44981.31  # L126
$ws.Cells.Item(126, 13).Value = -28980.251  # M126
$ws.Cells.Item(126, 14).Value = -49921.31  # N126

$ws.Cells.Item(132, 8).Value = 7939.136  # H132
$ws.Cells.Item(132, 9).Value = 8233.15  # I132
$ws.Cells.Item(132, 10).Value = 4999  # J132
$ws.Cells.Item(132, 11).Value = 24699.45  # K132
$ws.Cells.Item(132, 12).Value = 14997  # L132
$ws.Cells.Item(132, 13).Value = -22169.45  # M132
$ws.Cells.Item(132, 14).Value = -20057  # N132

$ws.Cells.Item(136, 8).Value = 6913.5806  # H136
$ws.Cells.Item(136, 9).Value = 7511.1333  # I136
$ws.Cells.Item(136, 10).Value = 6353.375  # J136
$ws.Cells.Item(136, 11).Value = 22533.3999  # K136
$ws.Cells.Item(136, 12).Value = 19060.125  # L136
$ws.Cells.Item(136, 13).Value = -19983.3999  # M136
$ws.Cells.Item(136, 14).Value = -24160.125  # N136

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 486.27274  # H113
$ws.Cells.Item(113, 10).Value = 829.2  # J113
$ws.Cells.Item(113, 12).Value = 2487.6  # L113
$ws.Cells.Item(113, 14).Value = -6827.6  # N113

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(94, 8).Value = 19999  # H94
$ws.Cells.Item(94, 10).Value = 19999  # J94
$ws.Cells.Item(94, 12).Value = 19999  # L94
$ws.Cells.Item(94, 14).Value = -21351  # N94

$ws.Cells.Item(99, 8).Value = 9482.799999999999  # H99
$ws.Cells.Item(99, 9).Value = 8606  # I99
$ws.Cells.Item(99, 10).Value = 12990  # J99
$ws.Cells.Item(99, 11).Value = 8606  # K99
$ws.Cells.Item(99, 12).Value = 12990  # L99
$ws.Cells.Item(99, 13).Value = -6360  # M99
$ws.Cells.Item(99, 14).Value = -17482  # N99

$ws.Cells.Item(102, 8).Value = 4218.8276  # H102
$ws.Cells.Item(102, 9).Value = 4239.409  # I102
$ws.Cells.Item(102, 11).Value = 4239.409  # K102
$ws.Cells.Item(102, 13).Value = -2617.409  # M102

$ws.Cells.Item(107, 8).Value = 530.6875  # H107
$ws.Cells.Item(107, 9).Value = 530.9231  # I107
$ws.Cells.Item(107, 11).Value = 530.9231  # K107
$ws.Cells.Item(107, 13).Value = 1389.0769  # M107

$ws.Cells.Item(113, 8).Value = 10761.77  # H113
$ws.Cells.Item(113, 9).Value = 3945.8572  # I113
$ws.Cells.Item(113, 10).Value = 18713.666  # J113
$ws.Cells.Item(113, 11).Value = 3945.8572  # K113
$ws.Cells.Item(113, 12).Value = 18713.666  # L113
$ws.Cells.Item(113, 13).Value = -1775.8572  # M113
$ws.Cells.Item(113, 14).Value = -23053.666  # N113

$ws.Cells.Item(126, 8).Value = 4647.087  # H126
$ws.Cells.Item(126, 9).Value = 4219.875  # I126
$ws.Cells.Item(126, 11).Value = 12659.625  # K126
$ws.Cells.Item(126, 13).Value = -10189.625  # M126

$ws.Cells.Item(132, 8).Value = 3310.75  # H132
$ws.Cells.Item(132, 9).Value = 2081.111  # I132
$ws.Cells.Item(132, 11).Value = 6243.333  # K132
$ws.Cells.Item(132, 13).Value = -3713.333  # M132

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(4, 8).Value = 2019  # H4
$ws.Cells.Item(4, 10).Value = 0  # J4
$ws.Cells.Item(4, 12).Value = 0  # L4
$ws.Cells.Item(4, 14).ClearContents()  # N4

$ws.Cells.Item(16, 8).Value = 3232.7  # H16
$ws.Cells.Item(16, 9).Value = 3232.7  # I16
$ws.Cells.Item(16, 10).Value = 0  # J16
$ws.Cells.Item(16, 11).Value = 3232.7  # K16
$ws.Cells.Item(16, 12).Value = 0  # L16
$ws.Cells.Item(16, 13).Value = -3062.7  # M16
$ws.Cells.Item(16, 14).ClearContents()  # N16

$ws.Cells.Item(28, 8).Value = 2019  # H28
$ws.Cells.Item(28, 10).Value = 0  # J28
$ws.Cells.Item(28, 12).Value = 0  # L28
$ws.Cells.Item(28, 14).ClearContents()  # N28

$ws.Cells.Item(37, 8).Value = 2019  # H37
$ws.Cells.Item(37, 10).Value = 0  # J37
$ws.Cells.Item(37, 12).Value = 0  # L37
$ws.Cells.Item(37, 14).ClearContents()  # N37

$ws.Cells.Item(41, 8).Value = 30000  # H41
$ws.Cells.Item(41, 9).Value = 30000  # I41
$ws.Cells.Item(41, 11).Value = 30000  # K41
$ws.Cells.Item(41, 13).Value = -29562  # M41

$ws.Cells.Item(50, 8).Value = 48976  # H50
$ws.Cells.Item(50, 9).Value = 48976  # I50
$ws.Cells.Item(50, 11).Value = 48976  # K50
$ws.Cells.Item(50, 13).Value = -48339  # M50

$ws.Cells.Item(122, 8).Value = 2969.1875  # H122
$ws.Cells.Item(122, 9).Value = 2175.5833  # I122
$ws.Cells.Item(122, 10).Value = 5350  # J122
$ws.Cells.Item(122, 11).Value = 6526.749899999999  # K122
$ws.Cells.Item(122, 12).Value = 16050  # L122
$ws.Cells.Item(122, 13).Value = -4076.749899999999  # M122
$ws.Cells.Item(122, 14).Value = -20950  # N122

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(3, 8).Value = 5773.6665  # H3
$ws.Cells.Item(3, 9).Value = 4660.5  # I3
$ws.Cells.Item(3, 10).Value = 8000  # J3
$ws.Cells.Item(3, 11).Value = 4660.5  # K3
$ws.Cells.Item(3, 12).Value = 8000  # L3
$ws.Cells.Item(3, 13).Value = -4546.5  # M3
$ws.Cells.Item(3, 14).Value = -8228  # N3

$ws.Cells.Item(34, 8).Value = 9899  # H34
$ws.Cells.Item(34, 9).Value = 9899  # I34
$ws.Cells.Item(34, 11).Value = 9899  # K34
$ws.Cells.Item(34, 13).Value = -9696  # M34

$ws.Cells.Item(37, 8).Value = 9899  # H37
$ws.Cells.Item(37, 9).Value = 9899  # I37
$ws.Cells.Item(37, 11).Value = 9899  # K37
$ws.Cells.Item(37, 13).Value = -9696  # M37

$ws.Cells.Item(40, 8).Value = 0  # H40
$ws.Cells.Item(40, 9).Value = 0  # I40
$ws.Cells.Item(40, 11).Value = 0  # K40
$ws.Cells.Item(40, 13).ClearContents()  # M40

$ws.Cells.Item(122, 8).Value = 5073.9556  # H122
$ws.Cells.Item(122, 9).Value = 3035.7715  # I122
$ws.Cells.Item(122, 10).Value = 12207.6  # J122
$ws.Cells.Item(122, 11).Value = 9107.3145  # K122
$ws.Cells.Item(122, 12).Value = 36622.8  # L122
$ws.Cells.Item(122, 13).Value = -6657.3145  # M122
$ws.Cells.Item(122, 14).Value = -41522.8  # N122

$ws.Cells.Item(136, 8).Value = 8962.210999999999  # H136
$ws.Cells.Item(136, 9).Value = 14641.714  # I136
$ws.Cells.Item(136, 11).Value = 43925.142  # K136
$ws.Cells.Item(136, 13).Value = -41375.142  # M136

